$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Underidentified CFA models ... have fewer free parameters ..." ->
#    "... have more free parameters ..."
#    (only this one paragraph - the sibling "Overidentified" / "Just-identified"
#    paragraphs re-use the same trailing sentence and must stay untouched)
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $ptext = $p.Range.Text
    if ($ptext -like "*Underidentified*fewer free parameters*") {

        $pRange = $p.Range
        $pRange.Find.Execute("fewer", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "more", 2)

        # Re-seat the "_GoBack" bookmark right after the word that was just
        # typed ("more"), matching where Word leaves the edit-point marker.
        $pStart = $p.Range.Start
        $newText = $p.Range.Text
        $offset = $newText.IndexOf(" free parameters")
        if ($offset -ge 0) {
            $target = $pStart + $offset
            $goBack = $d.Bookmarks.Item("_GoBack")
            $goBack.Delete()
            $targetRange = $d.Range($target, $target)
            $d.Bookmarks.Add("_GoBack", $targetRange)
        }
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Register the (until now unused) built-in "Balloon Text" style pair.
# ---------------------------------------------------------------------------
$balloon = $d.Styles.Add("BalloonText", 1)
$balloon.NameLocal = "Balloon Text"
$balloon.BaseStyle = $d.Styles.Item("Normal")
$balloon.Priority = 99
$balloon.UnhideWhenUsed = $true
$balloon.ParagraphFormat.SpaceAfter = 0
$balloon.ParagraphFormat.LineSpacingRule = 0
$balloon.Font.Name = "Segoe UI"
$balloon.Font.NameAscii = "Segoe UI"
$balloon.Font.NameBi = "Segoe UI"
$balloon.Font.Size = 9
$balloon.Font.SizeBi = 9

$balloonChar = $d.Styles.Add("BalloonTextChar", 2)
$balloonChar.NameLocal = "Balloon Text Char"
$balloonChar.BaseStyle = $d.Styles.Item("DefaultParagraphFont")
$balloonChar.Priority = 99
$balloonChar.Font.Name = "Segoe UI"
$balloonChar.Font.NameAscii = "Segoe UI"
$balloonChar.Font.NameBi = "Segoe UI"
$balloonChar.Font.Size = 9
$balloonChar.Font.SizeBi = 9

$balloon.LinkStyle = $balloonChar
$balloonChar.LinkStyle = $balloon
